# Apply updated TPM-derived metrics to Fgf2-Gpc4 LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.313179
$ws.Range("H2").Value = 0.939537
$ws.Range("I2").Value = 0.02707464596575709
$ws.Range("J2").Value = 0.0270746459657571
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 22.51161266666667
$ws.Range("N2").Value = 67.53483800000001
$ws.Range("O2").Value = 0.3150979864474181
$ws.Range("P2").Value = 0.3150979864474181
$ws.Range("Q2").Value = 7.050164343334
$ws.Range("R2").Value = 63.451479090006
$ws.Range("S2").Value = 0.008531166427586771
$ws.Range("T2").Value = 0.008531166427586773

$ws.Range("G3").Value = 0.313179
$ws.Range("H3").Value = 0.939537
$ws.Range("I3").Value = 0.02707464596575709
$ws.Range("J3").Value = 0.0270746459657571
$ws.Range("N3").Value = 90.745323
$ws.Range("O3").Value = 0.4233913844114141
$ws.Range("P3").Value = 0.4233913844114141
$ws.Range("Q3").Value = 9.473176503938999
$ws.Range("R3").Value = 85.258588535451
$ws.Range("S3").Value = 0.0114631718378908
$ws.Range("T3").Value = 0.0114631718378908

$ws.Range("G4").Value = 0.313179
$ws.Range("H4").Value = 0.939537
$ws.Range("I4").Value = 0.02707464596575709
$ws.Range("J4").Value = 0.0270746459657571
$ws.Range("M4").Value = 18.63333
$ws.Range("N4").Value = 55.89999
$ws.Range("O4").Value = 0.2608131567803688
$ws.Range("P4").Value = 0.2608131567803688
$ws.Range("Q4").Value = 5.83556765607
$ws.Range("R4").Value = 52.52010890463
$ws.Range("S4").Value = 0.007061423883039985
$ws.Range("T4").Value = 0.007061423883039985

$ws.Range("G5").Value = 0.313179
$ws.Range("H5").Value = 0.939537
$ws.Range("I5").Value = 0.02707464596575709
$ws.Range("J5").Value = 0.0270746459657571
$ws.Range("M5").Value = 0.04982966666666667
$ws.Range("N5").Value = 0.149489
$ws.Range("O5").Value = 0.000697472360799001
$ws.Range("P5").Value = 0.000697472360799001
$ws.Range("Q5").Value = 0.015605605177
$ws.Range("R5").Value = 0.140450446593
$ws.Range("S5").Value = 0.00001888381723953375
$ws.Range("T5").Value = 0.00001888381723953375

$ws.Range("I6").Value = 0.6982806158817221
$ws.Range("J6").Value = 0.6982806158817222
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 22.51161266666667
$ws.Range("N6").Value = 67.53483800000001
$ws.Range("O6").Value = 0.3150979864474181
$ws.Range("P6").Value = 0.3150979864474181
$ws.Range("Q6").Value = 181.8303776144302
$ws.Range("R6").Value = 1636.473398529872
$ws.Range("S6").Value = 0.2200268160395936
$ws.Range("T6").Value = 0.2200268160395936

$ws.Range("I7").Value = 0.6982806158817221
$ws.Range("J7").Value = 0.6982806158817222
$ws.Range("N7").Value = 90.745323
$ws.Range("O7").Value = 0.4233913844114141
$ws.Range("P7").Value = 0.4233913844114141
$ws.Range("S7").Value = 0.2956459966658172
$ws.Range("T7").Value = 0.2956459966658173

$ws.Range("I8").Value = 0.6982806158817221
$ws.Range("J8").Value = 0.6982806158817222
$ws.Range("M8").Value = 18.63333
$ws.Range("N8").Value = 55.89999
$ws.Range("O8").Value = 0.2608131567803688
$ws.Range("P8").Value = 0.2608131567803688
$ws.Range("Q8").Value = 150.50478525384
$ws.Range("R8").Value = 1354.54306728456
$ws.Range("S8").Value = 0.1821207717466521
$ws.Range("T8").Value = 0.1821207717466521

$ws.Range("I9").Value = 0.6982806158817221
$ws.Range("J9").Value = 0.6982806158817222
$ws.Range("M9").Value = 0.04982966666666667
$ws.Range("N9").Value = 0.149489
$ws.Range("O9").Value = 0.000697472360799001
$ws.Range("P9").Value = 0.000697472360799001
$ws.Range("Q9").Value = 0.4024832534462223
$ws.Range("R9").Value = 3.622349281016
$ws.Range("S9").Value = 0.0004870314296592051
$ws.Range("T9").Value = 0.0004870314296592052

$ws.Range("G10").Value = 2.897745666666667
$ws.Range("H10").Value = 8.693237
$ws.Range("I10").Value = 0.2505130868410934
$ws.Range("J10").Value = 0.2505130868410934
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.51161266666667
$ws.Range("N10").Value = 67.53483800000001
$ws.Range("O10").Value = 0.3150979864474181
$ws.Range("P10").Value = 0.3150979864474181
$ws.Range("Q10").Value = 65.23292805451179
$ws.Range("R10").Value = 587.096352490606
$ws.Range("S10").Value = 0.07893616924235571
$ws.Range("T10").Value = 0.07893616924235571

$ws.Range("G11").Value = 2.897745666666667
$ws.Range("H11").Value = 8.693237
$ws.Range("I11").Value = 0.2505130868410934
$ws.Range("J11").Value = 0.2505130868410934
$ws.Range("N11").Value = 90.745323
$ws.Range("O11").Value = 0.4233913844114141
$ws.Range("P11").Value = 0.4233913844114141
$ws.Range("Q11").Value = 87.65228883117234
$ws.Range("R11").Value = 788.8705994805509
$ws.Range("S11").Value = 0.1060650826508273
$ws.Range("T11").Value = 0.1060650826508273

$ws.Range("G12").Value = 2.897745666666667
$ws.Range("H12").Value = 8.693237
$ws.Range("I12").Value = 0.2505130868410934
$ws.Range("J12").Value = 0.2505130868410934
$ws.Range("M12").Value = 18.63333
$ws.Range("N12").Value = 55.89999
$ws.Range("O12").Value = 0.2608131567803688
$ws.Range("P12").Value = 0.2608131567803688
$ws.Range("Q12").Value = 53.99465126307
$ws.Range("R12").Value = 485.95186136763
$ws.Range("S12").Value = 0.06533710899382023
$ws.Range("T12").Value = 0.06533710899382021

$ws.Range("G13").Value = 2.897745666666667
$ws.Range("H13").Value = 8.693237
$ws.Range("I13").Value = 0.2505130868410934
$ws.Range("J13").Value = 0.2505130868410934
$ws.Range("M13").Value = 0.04982966666666667
$ws.Range("N13").Value = 0.149489
$ws.Range("O13").Value = 0.000697472360799001
$ws.Range("P13").Value = 0.000697472360799001
$ws.Range("Q13").Value = 0.1443937006547778
$ws.Range("R13").Value = 1.299543305893
$ws.Range("S13").Value = 0.0001747259540901025
$ws.Range("T13").Value = 0.0001747259540901025

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.2791366666666666
$ws.Range("H14").Value = 0.83741
$ws.Range("I14").Value = 0.02413165131142748
$ws.Range("J14").Value = 0.02413165131142749
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 22.51161266666667
$ws.Range("N14").Value = 67.53483800000001
$ws.Range("O14").Value = 0.3150979864474181
$ws.Range("P14").Value = 0.3150979864474181
$ws.Range("Q14").Value = 6.283816521064444
$ws.Range("R14").Value = 56.55434868958
$ws.Range("S14").Value = 0.007603834737881997
$ws.Range("T14").Value = 0.007603834737881999

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.2791366666666666
$ws.Range("H15").Value = 0.83741
$ws.Range("I15").Value = 0.02413165131142748
$ws.Range("J15").Value = 0.02413165131142749
$ws.Range("N15").Value = 90.745323
$ws.Range("O15").Value = 0.4233913844114141
$ws.Range("P15").Value = 0.4233913844114141
$ws.Range("Q15").Value = 8.443448992603333
$ws.Range("R15").Value = 75.99104093343
$ws.Range("S15").Value = 0.0102171332568788
$ws.Range("T15").Value = 0.0102171332568788

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.2791366666666666
$ws.Range("H16").Value = 0.83741
$ws.Range("I16").Value = 0.02413165131142748
$ws.Range("J16").Value = 0.02413165131142749
$ws.Range("M16").Value = 18.63333
$ws.Range("N16").Value = 55.89999
$ws.Range("O16").Value = 0.2608131567803688
$ws.Range("P16").Value = 0.2608131567803688
$ws.Range("Q16").Value = 5.201245625099999
$ws.Range("R16").Value = 46.8112106259
$ws.Range("S16").Value = 0.006293852156856529
$ws.Range("T16").Value = 0.006293852156856529

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.2791366666666666
$ws.Range("H17").Value = 0.83741
$ws.Range("I17").Value = 0.02413165131142748
$ws.Range("J17").Value = 0.02413165131142749
$ws.Range("M17").Value = 0.04982966666666667
$ws.Range("N17").Value = 0.149489
$ws.Range("O17").Value = 0.000697472360799001
$ws.Range("P17").Value = 0.000697472360799001
$ws.Range("Q17").Value = 0.01390928705444444
$ws.Range("R17").Value = 0.12518358349
$ws.Range("S17").Value = 0.00001683115981015964
$ws.Range("T17").Value = 0.00001683115981015964
